{"js": "const target = \"Literatura se v\\u00A0tomto obdob\u00ed d\u011blila na tradi\u010dn\u00ed a experiment\u00e1ln\u00ed\";\n\nconst results = context.document.body.search(target, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Proud Realistick\u00fd\", Word.InsertLocation.replace);\n} else {\n  // Fallback: NBSP normalized to a regular space, or other whitespace variance.\n  const fallback = target.replace(/\\u00A0/g, \" \");\n  const fallbackResults = context.document.body.search(fallback, { matchCase: true });\n  fallbackResults.load(\"items\");\n  await context.sync();\n  if (fallbackResults.items.length > 0) {\n    fallbackResults.items[0].insertText(\"Proud Realistick\u00fd\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$r.Find.Text = \"Literatura se v\" + [char]0x00A0 + \"tomto obdob\u00ed d\u011blila na tradi\u010dn\u00ed a experiment\u00e1ln\u00ed\"\n$r.Find.Replacement.Text = \"Proud Realistick\u00fd\"\n$r.Find.Execute($r.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $r.Find.Replacement.Text, 2)\n"}
